$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OrderPage")

# Update the two date cells that had the stray "6/19/2019" value so they
# match the other date entries ("6/18/2019").
$ws.Range("F4").Value = "6/18/2019"
$ws.Range("E5").Value = "6/18/2019"

# Update the active selection on the sheet.
$ws.Activate()
$ws.Range("E8").Select()
